$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36-52 down to 37-53
$ws.Rows(36).Insert()

# Populate the newly inserted row 36 with the new data record
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C36").Value = "Ñuble"
$ws.Range("D36").Value = 44830
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = 100112026
$ws.Range("G36").Value = "Haba"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 100
$ws.Range("K36").Value = 9000
$ws.Range("L36").Value = 9500
$ws.Range("M36").Value = 9250
$ws.Range("N36").Value = "$/saco 25 kilos"
$ws.Range("O36").Value = "Provincia de Limarí"
$ws.Range("P36").Value = 370
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
